$d = $word.ActiveDocument

# Heading3: Cellular Biology -> Cell Biology
$d.Content.Find.Execute("Cellular Biology", $true, $false, $false, $false, $false, $true, 1, $false, "Cell Biology", 2) | Out-Null

# Ativacao date: 2018 -> 2025
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# Objetivos (PT): remove "de biologia geral"
$d.Content.Find.Execute("Dotar os alunos dos conhecimentos de biologia celular abrangendo a organização estrutural e molecular da célula, proporcionando os conhecimentos básicos de biologia geral necessários à compreensão das demais disciplinas correlacionadas do curso de Engenharia Bioquímica.", $true, $false, $false, $false, $false, $true, 1, $false, "Dotar os alunos dos conhecimentos de biologia celular abrangendo a organização estrutural e molecular da célula, proporcionando os conhecimentos básicos necessários à compreensão das demais disciplinas correlacionadas do curso de Engenharia Bioquímica.", 2) | Out-Null

# Programa resumido (EN): rewrite short sentence
$d.Content.Find.Execute("Cellular origin and evolution; structural analysis of cells; internal organization of cells.", $true, $false, $false, $false, $false, $true, 1, $false, "Origin and evolution of cells; structural analysis of cells; internal organization of cells.", 2) | Out-Null

# Programa (PT): add leading "- " and normalize dash spacing
$d.Content.Find.Execute("Origem e evolução das células: Conceitos básicos de sistemática e filogenia molecular, características dos três domínios. –Análise estrutural das células ao microscópio: Microscopia ótica e microscopia eletrônica.–Organização interna das células: Células procarióticas e eucarióticas; estrutura e transporte através das membranas; compartimentos intracelulares (núcleo, retículo endoplasmático rugoso e liso, complexo de golgi, lisossomos e peroxissomos) e endereçamento de proteínas; tráfego intracelular de vesículas (via secretora e endocítica); conversão de energia (mitocôndria e cloroplasto); comunicação e sinalização celular; citoesqueleto; ciclo e divisão celular (mitose e meiose); matriz extracelulares e parede celular vegetal.", $true, $false, $false, $false, $false, $true, 1, $false, "- Origem e evolução das células: Conceitos básicos de sistemática e filogenia molecular, características dos três domínios. – Análise estrutural das células ao microscópio: Microscopia ótica e microscopia eletrônica. – Organização interna das células: Células procarióticas e eucarióticas; estrutura e transporte através das membranas; compartimentos intracelulares (núcleo, retículo endoplasmático rugoso e liso, complexo de golgi, lisossomos e peroxissomos) e endereçamento de proteínas; tráfego intracelular de vesículas (via secretora e endocítica); conversão de energia (mitocôndria e cloroplasto); comunicação e sinalização celular; citoesqueleto; ciclo e divisão celular (mitose e meiose); matriz extracelulares e parede celular vegetal.", 2) | Out-Null

# Programa (EN): full rewrite of long paragraph
$d.Content.Find.Execute("Origin and evolution of cells: basic concepts of systematic and molecular phylogeny, characteristics of the three domains. Structural analysis of cells at the microscope: optical and electronic microscopy.Internal organization of cells:Prokaryotic and eukaryotic cells; structure and transport through the membranes; intracelular compartments (nucleus, rough and smooth endoplasmic reticulum, golgi complex; lysosomes and peroxisomes) and protein addressing; intracellular traffic of vesicles (secretory and endocytic pathway); energy conversion (mitochondria and chloroplast); cellular communication and signalization; cytoskeleton; cellular cycle and division (mitosis and meiosis); extracellular matrix and vegetal cell wall.", $true, $false, $false, $false, $false, $true, 1, $false, "- Origin and evolution of cells: Basic concepts of molecular systematics and phylogeny, characteristics of the three domains.– Structural analysis of cells under a microscope: Optical microscopy and electron microscopy. – Internal organization of cells: Prokaryotic and eukaryotic cells; structure and transport across membranes; intracellular compartments (nucleus, rough and smooth endoplasmic reticulum, golgi complex, lysosomes and peroxisomes) and protein addressing; intracellular vesicle trafficking (secretory and endocytic pathway); energy conversion (mitochondria and chloroplast); cell communication and signaling; cytoskeleton; cell cycle and division (mitosis and meiosis); extracellular matrix and plant cell wall.", 2) | Out-Null

# Avaliacao: add space in MF formula text
$d.Content.Find.Execute("MF=Média finalMF = (P1 + P2) / 2", $true, $false, $false, $false, $false, $true, 1, $false, "MF=Média final MF = (P1 + P2) / 2", 2) | Out-Null

# Avaliacao: add space in Nota final formula text
$d.Content.Find.Execute("Nota final (NF)NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "Nota final (NF) NF = (MF + PR)/2, onde PR é uma prova de recuperação. Prova de recuperação (PR) para alunos com Média Final maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final igual ou maior do que 5,0.", 2) | Out-Null

# Bibliografia: full rewrite
$d.Content.Find.Execute("- Alberts, B., et al. Biologia Molecular da Célula, 5ed. Artmed Editora Ltda, 2010 - Cooper, G. M., Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3° Ed. 2007.- Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edição. 2010. - Lodish, H.; Berk, A.; Matsudaira, P.; Kaiser, C. A.; Krieger, M.; Scott, M. P.; Zipurky; Darnell. Biologia Celular e Molecular. 5ª Edição. Editora Artmed, 2005.- Raven, P. H., Evert, S. E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edição, 2016. - Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edição. 2009.- De Roberts, E.M.F.; Hibs, J. Bases da biologia celular e molecular. Editora Guanabara Koogan, 2006.", $true, $false, $false, $false, $false, $true, 1, $false, "Alberts, B., Bray, D., Hopkin, K., Johnson, A., Lewis, J., Raff, M., Roberts, K., Walter, P., Renard, A.E.B.A.G.   Fundamentos da Biologia Celular, 4ª edição, Ed. Artmed, 2017. ISBN 978-8582714058. Alberts, B., Johnson, A., Lewis, J., Morgan, D., Raff, M., Roberts, K., Walter, P., Wilson, J., Hunt, T. Biologia Molecular da Célula. 6a Edição, 2017. ISBN 978-8582714225. Cooper, G.M.; Robert, E.H. A célula: uma abordagem molecular. Artmed Editora Ltda, 3ª Edição, 2007. ISBN 978-8536308838. De Roberts, E.M.F., Hib, J. Bases da Biologia Celular e Molecular. Editora Guanabara Koogan, 16ª Edição, 2014. ISBN 978-8527723633. Junqueira e Carneiro. Biologia Celular e Molecular. 12ª Edição, Guanabara Koogan, 2023. ISBN 978-8527739337.  Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14 Edição, 2016. ISBN 978-8582712979. Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 8ª Edição,2014. ISBN 978-8527723626. Tortora, G.J.; Funke, B.R.; Case, C.L. Microbiologia. Editora Artmed, 12ª edição, 2016. ISBN 978-8582713532. Wasserman, S.A.; Monorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora, 12ª Edição, 2022. ISBN 978-6558820673.", 2) | Out-Null

# Objetivos (EN): the italic paragraph right after the Portuguese "Objetivos" paragraph
# (which ends in "Engenharia Bioquímica.") is empty; insert its English translation
# into that existing (italic) run rather than hardcoding a paragraph index.
$objIndex = -1
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $cur = $d.Paragraphs.Item($i)
    if ($cur.Range.Text -like "*Engenharia Bioqu*") {
        $objIndex = $i + 1
        break
    }
}
if ($objIndex -ne -1) {
    $d.Paragraphs.Item($objIndex).Range.InsertAfter("Provide students with knowledge of cell biology, covering the structural and molecular organization of the cell, providing the basic knowledge necessary to understand the other related disciplines of the Biochemical Engineering course.")
}

